$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.485.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.934.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2878"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06646"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "108.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +27.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.918.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07613"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.177"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.6623"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "307.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +24.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "30.499.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007598"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.176.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.300"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.322"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.326"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.057"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1110"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.363"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.100"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.952"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05029"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7448"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.746"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01966"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.690"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.045"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8817"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.814"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4204"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.277"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.238"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1216"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05628"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.40%  "
